{"js": "// Title changes from:\n//   \"Chill Portions and yield using actual data\"\n// to:\n//   \"Estimating yield as a function of chill accumulation\"\n//\n// We operate only on the first paragraph (the document Title) so that\n// common words (\"data\", \"actual\", \"using\", ...) that also occur later\n// in the abstract body text are left untouched.\n\nconst body = context.document.body;\nbody.load(\"paragraphs\");\nawait context.sync();\n\nconst titlePara = body.paragraphs.getFirst();\n\n// Replace a single word/phrase found inside the title paragraph with new\n// text. Operating on a search hit confined to the title paragraph's own\n// range keeps the rest of the document safe from accidental matches, and\n// reusing the matched run via \"Replace\" keeps a clean one-run-per-token\n// structure identical to the original document.\nasync function replaceInTitle(searchText, replacement) {\n  const titleRange = titlePara.getRange();\n  const results = titleRange.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Chill\" -> \"Estimating\"\nawait replaceInTitle(\"Chill\", \"Estimating\");\n\n// 2) Drop \" Portions and\" (keeps the following space that precedes\n//    \"yield\" untouched, so \"Estimating yield ...\" reads correctly).\nawait replaceInTitle(\" Portions and\", \"\");\n\n// 3) \"using\" -> \"as\"\nawait replaceInTitle(\"using\", \"as\");\n\n// 4) \"actual\" -> \"a\"\nawait replaceInTitle(\"actual\", \"a\");\n\n// 5) \"data\" -> \"function\"\nawait replaceInTitle(\"data\", \"function\");\n\n// 6) Append \" of chill accumulation\" as individual word/space runs (to\n//    match the original document's one-run-per-word/space convention)\n//    right after the word \"function\" at the end of the title paragraph.\nfunction runXml(text) {\n  return '<w:r><w:t xml:space=\"preserve\">' + text + \"</w:t></w:r>\";\n}\n\nconst newTokens = [\" \", \"of\", \" \", \"chill\", \" \", \"accumulation\"];\nconst runsXml = newTokens.map(runXml).join(\"\");\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p>\" +\n  runsXml +\n  \"</w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nconst endRange = titlePara.getRange(\"End\");\nendRange.insertOoxml(ooxml, \"End\");\nawait context.sync();\n", "ps1": "# Title changes from:\n#   \"Chill Portions and yield using actual data\"\n# to:\n#   \"Estimating yield as a function of chill accumulation\"\n#\n# We work only against the first paragraph (the document Title) and\n# rebuild it run-by-run (one run per word / per single space, exactly\n# like the original paragraph), so paragraphs elsewhere in the abstract\n# that happen to share words (\"data\", \"actual\", \"using\", ...) are left\n# completely untouched.\n\n$d = $word.ActiveDocument\n$para = $d.Paragraphs(1)\n$styleName = $para.Style.NameLocal\n\n# Paragraph text without the trailing paragraph-mark character.\n$text = $para.Range.Text\n$text = $text.TrimEnd([char]13, [char]7)\n\n# Tokenize into words and single spaces, alternating - this mirrors the\n# paragraph's existing one-run-per-word/space convention.\n$tokens = @()\n$cur = \"\"\nforeach ($ch in $text.ToCharArray()) {\n    if ($ch -eq ' ') {\n        if ($cur -ne \"\") { $tokens += $cur; $cur = \"\" }\n        $tokens += \" \"\n    } else {\n        $cur += $ch\n    }\n}\nif ($cur -ne \"\") { $tokens += $cur }\n\n# Drop \"Portions\" and \"and\" (together with the spacing around them) -\n# i.e. remove the four-token run \" Portions and\" that follows \"Chill\",\n# leaving the single space that precedes \"yield\" intact.\n$portionsIdx = -1\nfor ($i = 0; $i -lt $tokens.Count; $i++) {\n    if ($tokens[$i] -eq \"Portions\") { $portionsIdx = $i; break }\n}\nif ($portionsIdx -ge 1 -and $tokens[$portionsIdx - 1] -eq \" \" -and\n    $tokens[$portionsIdx + 1] -eq \" \" -and $tokens[$portionsIdx + 2] -eq \"and\") {\n    $newTokens = @()\n    for ($i = 0; $i -lt $tokens.Count; $i++) {\n        if ($i -ge ($portionsIdx - 1) -and $i -le ($portionsIdx + 2)) { continue }\n        $newTokens += $tokens[$i]\n    }\n    $tokens = $newTokens\n}\n\n# Word-for-word renames.\nfor ($i = 0; $i -lt $tokens.Count; $i++) {\n    if ($tokens[$i] -eq \"Chill\") { $tokens[$i] = \"Estimating\" }\n    elseif ($tokens[$i] -eq \"using\") { $tokens[$i] = \"as\" }\n    elseif ($tokens[$i] -eq \"actual\") { $tokens[$i] = \"a\" }\n    elseif ($tokens[$i] -eq \"data\") { $tokens[$i] = \"function\" }\n}\n\n# Append \" of chill accumulation\" (as individual word/space tokens)\n# right after \"function\" at the end of the title.\n$funcIdx = -1\nfor ($i = 0; $i -lt $tokens.Count; $i++) {\n    if ($tokens[$i] -eq \"function\") { $funcIdx = $i; break }\n}\n\n$newTail = @(\" \", \"of\", \" \", \"chill\", \" \", \"accumulation\")\n$finalTokens = @()\nfor ($i = 0; $i -lt $tokens.Count; $i++) {\n    $finalTokens += $tokens[$i]\n    if ($i -eq $funcIdx) {\n        foreach ($tok in $newTail) { $finalTokens += $tok }\n    }\n}\n\n# Rebuild the paragraph as one <w:r> per token (word or space), keeping\n# its original paragraph style, and push it back in via WordOpenXML so\n# the run-per-token layout matches the source document's convention.\n$runsXml = \"\"\nforeach ($tok in $finalTokens) {\n    $escaped = $tok.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n    $runsXml += '<w:r><w:t xml:space=\"preserve\">' + $escaped + '</w:t></w:r>'\n}\n\n$paraXml = '<w:p><w:pPr><w:pStyle w:val=\"' + $styleName + '\"/></w:pPr>' + $runsXml + '</w:p>'\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $paraXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$para.Range.InsertXML($xml)\n"}
